# Add a new worksheet named "rho" at the end of the workbook (after "x matrix")
# and populate it with the rho matrix data, mirroring the structure used by
# the existing "A matrix" / "r matrix" / "x matrix" sheets.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "rho"

# A1 label
$ws.Cells.Item(1, 1).Value = "rho"

# Column headers B1:AV1 = 1..47
$headerRow = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47)
for ($j = 0; $j -lt $headerRow.Length; $j++) {
    $ws.Cells.Item(1, $j + 2).Value = $headerRow[$j]
}

# Row labels A2:A5 = 1..4
$rowLabels = @(1,2,3,4)
for ($i = 0; $i -lt $rowLabels.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $rowLabels[$i]
}

# Data block B2:AV5
$matrix = @(
    @(2,3,2,3,4,5,6,7,8,11,6,11,10,4,14,17,16,3,16,19,16,21,22,16,2,25,17,26,26,2,0,10,19,20,22,23,25,29,1,2,10,19,20,22,23,25,29),
    @(30,25,4,5,6,7,8,9,39,13,10,13,12,13,16,19,18,17,20,34,22,23,24,23,26,27,26,29,28,40,0,41,42,43,44,45,46,47,9,30,32,33,34,35,36,37,38),
    @(40,30,18,14,0,11,0,0,0,32,12,0,14,15,0,21,27,0,33,43,0,35,36,0,37,28,0,0,38,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,40,0,0,0,31,0,0,0,41,0,0,0,0,0,24,0,0,42,0,0,44,45,0,46,29,0,0,47,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
)
for ($i = 0; $i -lt $matrix.Length; $i++) {
    $row = $matrix[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 2).Value = $row[$j]
    }
}

# Column widths to match the other matrix-style sheets
$ws.Range("B1:AW1").ColumnWidth = 3

# Selection on the new sheet
$ws.Range("A6").Select()

Write-Output "rho sheet added"
